# Fruta / hortaliza, semanal
# Insert two new price records for "Doctor Davis" durazno (row 263-264),
# pushing the existing rows 263..351 down to 265..353.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 263-264; everything below shifts down by 2 rows.
$ws.Rows("263:264").Insert()

# ---- New row 263: Doctor Davis / Especial ----
$ws.Cells.Item(263, 1).Value = 5
$ws.Cells.Item(263, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(263, 3).Value = "Maule"
$ws.Cells.Item(263, 4).Value = 44609
$ws.Cells.Item(263, 5).Value = 7
$ws.Cells.Item(263, 6).Value = "Fruta"
$ws.Cells.Item(263, 7).Value = 100103
$ws.Cells.Item(263, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(263, 9).Value = 100103004
$ws.Cells.Item(263, 10).Value = "Durazno"
$ws.Cells.Item(263, 11).Value = "Doctor Davis"
$ws.Cells.Item(263, 12).Value = "Especial"
$ws.Cells.Item(263, 13).Value = 300
$ws.Cells.Item(263, 14).Value = 13000
$ws.Cells.Item(263, 15).Value = 13000
$ws.Cells.Item(263, 16).Value = 13000
$ws.Cells.Item(263, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(263, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(263, 19).Value = 867
$ws.Cells.Item(263, 20).Value = 15

# ---- New row 264: Doctor Davis / Extra (doble especial) ----
$ws.Cells.Item(264, 1).Value = 5
$ws.Cells.Item(264, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(264, 3).Value = "Maule"
$ws.Cells.Item(264, 4).Value = 44609
$ws.Cells.Item(264, 5).Value = 7
$ws.Cells.Item(264, 6).Value = "Fruta"
$ws.Cells.Item(264, 7).Value = 100103
$ws.Cells.Item(264, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(264, 9).Value = 100103004
$ws.Cells.Item(264, 10).Value = "Durazno"
$ws.Cells.Item(264, 11).Value = "Doctor Davis"
$ws.Cells.Item(264, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(264, 13).Value = 300
$ws.Cells.Item(264, 14).Value = 15000
$ws.Cells.Item(264, 15).Value = 15000
$ws.Cells.Item(264, 16).Value = 15000
$ws.Cells.Item(264, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(264, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(264, 19).Value = 1000
$ws.Cells.Item(264, 20).Value = 15
